# Apply the MIPS control-mapping fix:
#  1. Column I (I3:I18) checked for "sw" instead of "lw" - fix the opcode test.
#  2. Update the helper note text in E20.
#  3. Move the active selection to E20 (where the edit ended up / was last reviewed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the I-column formulas: they should test for "lw", matching the K column,
#    not "sw" (which is already tested by the L column).
$ws.Range("I3:I18").FormulaR1C1 = '=IF(RC[-4]="lw", 1, 0)'

# 2. Update the group-note text (shared string) shown in E20.
$ws.Range("E20").Value = "update this column if series changes"

# 3. Leave the cell selection on E20, matching the reviewed/active cell.
$ws.Range("E20").Select() | Out-Null
